# Update the "Inferences" list on Sheet1: the old 3-bullet list (rows 12-14)
# is replaced with a new, more detailed 5-bullet list (rows 12-16). The
# "Inferences:" heading itself (row 11) keeps its text/style unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "1. Maize(corn) : Can be grown using less volume of water, precipitation has negative impact on corn yield."
$ws.Range("A13").Value = "2. Rice: Gross enrolment ratio is a significant positive predictor."
$ws.Range("A14").Value = "3. Sugar Cane: Precipitation and water use efficiency  are major negative predictors which means sugarcane requires more water to grow."
$ws.Range("A15").Value = "4. Bananas: Requires high amount of government expenditure as Agriculture share of Government Expenditure is a significant positive predictor."
$ws.Range("A16").Value = "5. Sweet Potatoes: Precipitation has negative impact on its yield. Gross enrolment ratio is a significant positive predictor.  "

# Match the author's last-saved selection.
$null = $ws.Range("E20").Select()
